# feat: add 2022-Q4 data
#
# Summary of the change being applied:
#  - The old "2022-Q3" data sheet is duplicated so its original content is
#    preserved under a (new) "2022-Q3" tab placed after it.
#  - The original "2022-Q3" sheet is renamed to "2022-Q4" and its values are
#    updated to the new quarter's figures.
#  - The "总计" (summary) sheet gets its existing "2022-Q3" row updated to
#    "2022-Q4" (with the new total holding value) and a new row is appended
#    below it restoring the original "2022-Q3" summary figures.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)   # "总计"
$q3 = $wb.Worksheets.Item(2)        # "2022-Q3" (will become "2022-Q4")

# 1) Duplicate the "2022-Q3" sheet; the copy keeps the original data/format
#    and is placed right after the source sheet.
$q3.Copy($null, $q3)
$q3copy = $wb.Worksheets.Item(3)

# 2) Turn the original sheet into the "2022-Q4" sheet with its new figures,
#    then rename the duplicate back to "2022-Q3" (done in this order to avoid
#    a transient duplicate-name clash).
$q3.Name = "2022-Q4"
$q3copy.Name = "2022-Q3"

$q3.Range("D2").Value = 1.52
$q3.Range("E2").Value = 92.93
$q3.Range("F2").Value = 3.87
$q3.Range("G2").Value = 0.0588
$q3.Range("H2").Value = 9

# 3) Update the "总计" sheet: existing row becomes the "2022-Q4" entry, and a
#    new row is added below with the original "2022-Q3" entry.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("D2").Value = 0.06

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.05
